# Ghép chức năng và giao diện cho phần thông báo, quy trình, ban đào tạo
#
# Adds a new "Nội dung" (Content) column (F) to Sheet1, with one value per
# existing data row (rows 2-9), and moves the active selection to I3.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column header
$ws.Range("F1").Value = "Nội dung"

# New column values for the existing 8 data rows
$ws.Range("F2").Value = "àhajd"
$ws.Range("F3").Value = "fjadshfjdas"
$ws.Range("F4").Value = "ádjfdfsak"
$ws.Range("F5").Value = "ạdkahdjf"
$ws.Range("F6").Value = "ạkfdahsfja"
$ws.Range("F7").Value = "kdfkdhf"
$ws.Range("F8").Value = "dkfhkfd"
$ws.Range("F9").Value = "kdjfdjh"

# Update the active selection shown when the workbook is reopened
$ws.Range("I3").Select()
